$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously jumped from row 1 straight to row 3 (row 2 was
# missing). Add the missing row 2 with a single-space value in A2 - this
# appends a new shared string (" ") without disturbing any existing rows.
$ws.Range("A2").Value = " "
